# Append two new feed-log rows (rows 50 and 51) to Sheet1, following on from
# the existing data which currently ends at row 49 (A1:E49).
#
# Columns: A=run_id, B=rss_url_id, C=date, D=response, E=item_count

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last populated row in column A so the new rows are appended
# immediately after the existing data, regardless of current sheet state.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow1 = $lastRow + 1
$newRow2 = $lastRow + 2

# Row 50: run_id 49
$ws.Cells.Item($newRow1, 1).Value = 49
$ws.Cells.Item($newRow1, 2).Value = 1
$ws.Cells.Item($newRow1, 3).Value = "2024-06-15 22:12:10"
$ws.Cells.Item($newRow1, 4).Value = 200
$ws.Cells.Item($newRow1, 5).Value = 2

# Row 51: run_id 50
$ws.Cells.Item($newRow2, 1).Value = 50
$ws.Cells.Item($newRow2, 2).Value = 2
$ws.Cells.Item($newRow2, 3).Value = "2024-06-15 22:12:11"
$ws.Cells.Item($newRow2, 4).Value = 200
$ws.Cells.Item($newRow2, 5).Value = 0
